# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Rewrites the worker/period detail table (rows 16-40) with the refreshed
# data set: existing workers re-grouped together, new workers (CRISTIAN
# ROCHA JIMENEZ, SEBASTIAN ENRIQUE VILLARREAL PERTUZ) interleaved, and new
# period rows (2109-2205) appended per worker, plus updated "Valor Mora"
# (column G) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: C = N Doc Trabajador, D = Nombre Trabajador, E = Periodo Mora,
#          F = Valor Mora (date-serial looking number), G = Salario Basico
$rows = @(
    @{ R=16; C="1052966060"; D="RODOLFO VALENZUELA GONZALEZ";         E="2105"; F=36341; G=880000 },
    @{ R=17; C="19752336";   D="CRISTIAN ROCHA JIMENEZ";               E="2105"; F=36341; G=880000 },
    @{ R=18; C="1143377371"; D="SEBASTIAN ENRIQUE VILLARREAL PERTUZ";  E="2109"; F=36341; G=908526 },
    @{ R=19; C="1143377371"; D="SEBASTIAN ENRIQUE VILLARREAL PERTUZ";  E="2110"; F=36341; G=908526 },
    @{ R=20; C="1131106099"; D="DAVINSON DARIO DIAZ LEDESMA";          E="2111"; F=36341; G=877803 },
    @{ R=21; C="73181246";   D="MARCO AURELIO MORA GUZMAN";            E="2111"; F=35112; G=877803 },
    @{ R=22; C="1143377371"; D="SEBASTIAN ENRIQUE VILLARREAL PERTUZ";  E="2111"; F=36341; G=908526 },
    @{ R=23; C="1131106099"; D="DAVINSON DARIO DIAZ LEDESMA";          E="2112"; F=36341; G=877803 },
    @{ R=24; C="73181246";   D="MARCO AURELIO MORA GUZMAN";            E="2112"; F=35112; G=877803 },
    @{ R=25; C="1143377371"; D="SEBASTIAN ENRIQUE VILLARREAL PERTUZ";  E="2112"; F=36341; G=908526 },
    @{ R=26; C="1131106099"; D="DAVINSON DARIO DIAZ LEDESMA";          E="2201"; F=36341; G=877803 },
    @{ R=27; C="73181246";   D="MARCO AURELIO MORA GUZMAN";            E="2201"; F=35112; G=877803 },
    @{ R=28; C="1143377371"; D="SEBASTIAN ENRIQUE VILLARREAL PERTUZ";  E="2201"; F=36341; G=908526 },
    @{ R=29; C="1131106099"; D="DAVINSON DARIO DIAZ LEDESMA";          E="2202"; F=36341; G=877803 },
    @{ R=30; C="73181246";   D="MARCO AURELIO MORA GUZMAN";            E="2202"; F=35112; G=877803 },
    @{ R=31; C="1143377371"; D="SEBASTIAN ENRIQUE VILLARREAL PERTUZ";  E="2202"; F=36341; G=908526 },
    @{ R=32; C="1131106099"; D="DAVINSON DARIO DIAZ LEDESMA";          E="2203"; F=36341; G=877803 },
    @{ R=33; C="73181246";   D="MARCO AURELIO MORA GUZMAN";            E="2203"; F=35112; G=877803 },
    @{ R=34; C="1143377371"; D="SEBASTIAN ENRIQUE VILLARREAL PERTUZ";  E="2203"; F=36341; G=908526 },
    @{ R=35; C="1131106099"; D="DAVINSON DARIO DIAZ LEDESMA";          E="2204"; F=35112; G=877803 },
    @{ R=36; C="73181246";   D="MARCO AURELIO MORA GUZMAN";            E="2204"; F=35112; G=877803 },
    @{ R=37; C="1143377371"; D="SEBASTIAN ENRIQUE VILLARREAL PERTUZ";  E="2204"; F=36341; G=908526 },
    @{ R=38; C="1131106099"; D="DAVINSON DARIO DIAZ LEDESMA";          E="2205"; F=26919; G=877803 },
    @{ R=39; C="73181246";   D="MARCO AURELIO MORA GUZMAN";            E="2205"; F=26919; G=877803 },
    @{ R=40; C="1143377371"; D="SEBASTIAN ENRIQUE VILLARREAL PERTUZ";  E="2205"; F=27861; G=908526 }
)

# Clear the whole detail block first so the shared-strings table is rebuilt
# cleanly (matching the freshly re-sorted/re-grouped data) instead of
# retaining stale index positions from the previous row order.
$ws.Range("C16:G40").ClearContents()

foreach ($row in $rows) {
    $r = $row.R
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
}
